$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural shift -------------------------------------------------
# Insert a new blank row at the top; every existing row moves down by one
# (old row 1 -> row 2, ... old row 10 -> row 11).
$ws.Rows("1:1").Insert()

# Insert two new blank columns before the old "D" column; old D,E,F -> F,G,H.
$ws.Columns("D:E").Insert()

# --- Column widths ------------------------------------------------------
$ws.Columns("C:C").ColumnWidth = 13
$ws.Columns("D:E").ColumnWidth = 13

# --- New currency number format -----------------------------------------
$curFmt = '"£"#,##0.00;[Red]\-"£"#,##0.00'

# --- Row 2 (Stepper Motor NEMA 17) --------------------------------------
$ws.Range("B2").Value = 2
$ws.Range("C2").NumberFormat = $curFmt
$ws.Range("C2").Value = 12.21
$ws.Range("D2").NumberFormat = $curFmt
$ws.Range("D2").Value = 14.651999999999999
$ws.Range("E2").NumberFormat = $curFmt
$ws.Range("E2").Formula = "=B2*D2"

# --- Row 11 (Electromagnet) updated to the 25kg Farnell part ------------
$farnellUrl = "https://uk.farnell.com/dfrobot/dfr0800/electromagnetic-lock-5v-25kg/dp/3769906?gclid=Cj0KCQiAm5ycBhCXARIsAPldzoXpMsX7goYBt_3BHIXdxOR6pkPWgVRATAw-KLtvC3f8cBxiMN1z8-caApd1EALw_wcB&mckv=_dc|pcrid||plid||kword||match||slid||product|3769906|pgrid||ptaid|&CMP=KNC-GUK-GEN-SHOPPING-SMART-PMAX-Medium_ROAS-Test958&gross_price=true"

$ws.Range("A11").Value = "Electromagnet 5V 25kg"
$ws.Range("D11").NumberFormat = $curFmt
$ws.Range("D11").Value = 10.54
$ws.Range("E11").NumberFormat = $curFmt
$ws.Range("E11").Value = 10.54
$ws.Range("G11").Value = "Farnell"
$ws.Range("H11").Value = $farnellUrl

# --- Hyperlinks: the row/column Insert() calls above leave the two
#     existing hyperlinks pointing at their stale (pre-shift) refs, so
#     drop them all and re-add at the correct post-shift locations. ------
$ws.Hyperlinks.Delete()

$digikeyUrl = "https://www.digikey.co.uk/en/products/detail/adafruit-industries-llc/324/5022791?utm_adgroup=General&utm_source=google&utm_medium=cpc&utm_campaign=PMax:%20Smart%20Shopping_Product_Zombie%20SKUs&utm_term=&productid=5022791&gclid=CjwKCAiA7IGcBhA8EiwAFfUDsTNRLcWA6HjCHboznJs4vu9jdAmZ9ACY08ebfRRH66cpmNJiw-DjmxoCTKYQAvD_BwE"
$ws.Hyperlinks.Add($ws.Range("H2"), $digikeyUrl, "", "", $digikeyUrl)

$mouserUrl = "https://www.mouser.co.uk/ProductDetail/Espressif-Systems/ESP32-S2-DevKitC-1RU?qs=pBJMDPsKWf1wWYxkgQUBfQ%3D%3D&mgh=1&vip=1&gclid=CjwKCAiA7IGcBhA8EiwAFfUDsYRj10h3izk5vQ4UWmTe7Dn7_zCRgf5sXVQpd2JwSykeqy0EkNWQARoC7jYQAvD_BwE"
$ws.Hyperlinks.Add($ws.Range("H4"), $mouserUrl)

$ws.Hyperlinks.Add($ws.Range("H11"), $farnellUrl, "", "", $farnellUrl)

# --- Page setup (paper size / orientation) -------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection matches the saved file ------------------------------------
$ws.Range("C2").Select()
